$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; existing rows 12:75 shift down to 13:76.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new data record.
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "Femacal de La Calera"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44901
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 100112022
$ws.Range("G12").Value = "Arveja Verde"
$ws.Range("H12").Value = "Perfection"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 65
$ws.Range("K12").Value = 18000
$ws.Range("L12").Value = 19000
$ws.Range("M12").Value = 18462
$ws.Range("N12").Value = "$/saco 25 kilos"
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 738
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
